# The source data set gained one more weekly record (week of 44776),
# which pushes every existing "Apio" row at La Palmera de La Serena down
# by two rows (one row per "Primera"/"Segunda" quality pair).
#
# Insert two blank rows right before the current row 409 so everything
# that used to live at 409..452 now lives at 411..454, then populate the
# two freshly inserted rows (409/410) with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(409).Insert()
$ws.Rows.Item(409).Insert()

# New row 409 - "Primera" quality
$ws.Range("A409").Value = 8
$ws.Range("B409").Value = "Terminal La Palmera de La Serena"
$ws.Range("C409").Value = "Coquimbo"
$ws.Range("D409").Value = 44776
$ws.Range("E409").Value = 4
$ws.Range("F409").Value = 100112017
$ws.Range("G409").Value = "Apio"
$ws.Range("H409").Value = "Americana (o)"
$ws.Range("I409").Value = "Primera"
$ws.Range("J409").Value = 2460
$ws.Range("K409").Value = 8000
$ws.Range("L409").Value = 9000
$ws.Range("M409").Value = 8500
$ws.Range("N409").Value = "$/docena de matas"
$ws.Range("O409").Value = "Provincia del Elquí"
$ws.Range("P409").Value = 1417
$ws.Range("Q409").Value = 6
$ws.Range("R409").Value = "Hortaliza"

# New row 410 - "Segunda" quality
$ws.Range("A410").Value = 8
$ws.Range("B410").Value = "Terminal La Palmera de La Serena"
$ws.Range("C410").Value = "Coquimbo"
$ws.Range("D410").Value = 44776
$ws.Range("E410").Value = 4
$ws.Range("F410").Value = 100112017
$ws.Range("G410").Value = "Apio"
$ws.Range("H410").Value = "Americana (o)"
$ws.Range("I410").Value = "Segunda"
$ws.Range("J410").Value = 1320
$ws.Range("K410").Value = 6500
$ws.Range("L410").Value = 7000
$ws.Range("M410").Value = 6750
$ws.Range("N410").Value = "$/docena de matas"
$ws.Range("O410").Value = "Provincia del Elquí"
$ws.Range("P410").Value = 1125
$ws.Range("Q410").Value = 6
$ws.Range("R410").Value = "Hortaliza"

Write-Output "Inserted rows 409-410 and shifted remaining data to 411-454"
